$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (H) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 440
$wsOff.Range("C2").Value = 304
$wsOff.Range("D2").Value = 116
$wsOff.Range("E2").Value = 54

# Sheet "DEF" - row 2 (H) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 539
$wsDef.Range("C2").Value = 375
$wsDef.Range("D2").Value = 109
$wsDef.Range("E2").Value = 52
$wsDef.Range("F2").Value = 12
